$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75, pushing existing rows 75-86 down to 76-87
$ws.Rows("75:75").Insert()

# Populate the newly inserted row 75 with the new weekly data point
$ws.Cells.Item(75, 1).Value = 8
$ws.Cells.Item(75, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(75, 3).Value = "Coquimbo"
$ws.Cells.Item(75, 4).Value = 44476
$ws.Cells.Item(75, 5).Value = 4
$ws.Cells.Item(75, 6).Value = 100112044
$ws.Cells.Item(75, 7).Value = "Perejil"
$ws.Cells.Item(75, 8).Value = "Sin especificar"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 3000
$ws.Cells.Item(75, 11).Value = 1500
$ws.Cells.Item(75, 12).Value = 2000
$ws.Cells.Item(75, 13).Value = 1750
$ws.Cells.Item(75, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(75, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(75, 16).Value = 1167
$ws.Cells.Item(75, 17).Value = 1.5
$ws.Cells.Item(75, 18).Value = "Hortaliza"
